$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Ready %" column header (cell G3, part of the "Tasks" table) to "Ready pct"
$ws.Range("G3").Value = "Ready pct"

# Restore the cursor/selection left behind by the edit (cell I3)
$ws.Range("I3").Select()
